# "Generate Report for Handoff" — localization status moves from
# "In Translation" to "Ready for handoff", and the handoff/report
# timestamps on the Overview, zh-cn and de-de sheets advance to the
# moment the report was (re)generated.

$wb = $excel.ActiveWorkbook

$oldStatus = "In Translation"
$newStatus = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus              # zh-cn status
$wsOverview.Range("F2").Value = $newStatus              # de-de status
$wsOverview.Range("G2").Value = "2016-08-13 01:13:09"    # Latest HO Xliff Generate Date

# Status text got longer ("In Translation" -> "Ready for handoff"), so the
# Status columns re-autofit to a wider column width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus                   # Status
$wsZhCn.Range("H2").Value = "2016-08-13 01:12:59"         # Latest Handoff Datetime
$wsZhCn.Columns.Item(3).ColumnWidth = 16.38265482584637

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus                    # Status
$wsDeDe.Range("H2").Value = "2016-08-13 01:13:09"          # Latest Handoff Datetime
$wsDeDe.Columns.Item(3).ColumnWidth = 16.38265482584637
